# Apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores figures as plain text (it uses "." as a thousands
# separator, e.g. "37.261.18"). Force text formatting before writing any D-column
# value that would otherwise be auto-parsed as a plain number (single "." decimal),
# so it round-trips the same way the original inline strings did.

$ws.Range("D2").Value = "37.261.18"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.028.19"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.52"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.37"
$ws.Range("E8").Value = "  -4.19%  "
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  -4.96%  "
$ws.Range("D12").Value = "2.327.14"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.29"
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.39"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.745"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "2.015.75"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "37.229.58"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.27"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.68"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  -5.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.32"
$ws.Range("E26").Value = "  -6.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.48"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.128"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.83"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("E35").Value = "  -5.77%  "
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.48"
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0218"
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").Value = "1.481.08"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.19"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0923"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.44"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.77"
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.23"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "2.214.29"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.63"
$ws.Range("E51").Value = "  -9.70%  "
